$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1299
$ws.Range("I28").Value = 667.6667
$ws.Range("J28").Value = 5087
$ws.Range("K28").Value = 667.6667
$ws.Range("L28").Value = 5087
$ws.Range("M28").Value = -182.6667
$ws.Range("N28").Value = -6057
$ws.Range("H51").Value = 13894506
$ws.Range("J51").Value = 6759
$ws.Range("L51").Value = 6759
$ws.Range("N51").Value = -7727
$ws.Range("H129").Value = 1241.1111
$ws.Range("I129").Value = 1081.7142
$ws.Range("J129").Value = 1799
$ws.Range("K129").Value = 3245.1426
$ws.Range("L129").Value = 5397
$ws.Range("M129").Value = 1754.8574
$ws.Range("N129").Value = -15397
$ws.Range("H138").Value = 2619.6667
$ws.Range("I138").Value = 2082.762
$ws.Range("J138").Value = 3872.4443
$ws.Range("K138").Value = 6248.286
$ws.Range("L138").Value = 11617.3329
$ws.Range("M138").Value = -1108.286
$ws.Range("N138").Value = -21897.3329

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17612.486
$ws.Range("I32").Value = 17834
$ws.Range("K32").Value = 17834
$ws.Range("M32").Value = -17547
$ws.Range("H122").Value = 2532.2
$ws.Range("J122").Value = 4133
$ws.Range("L122").Value = 12399
$ws.Range("N122").Value = -17299
$ws.Range("H132").Value = 1698.4073
$ws.Range("I132").Value = 1643.579
$ws.Range("J132").Value = 1828.625
$ws.Range("K132").Value = 4930.737
$ws.Range("L132").Value = 5485.875
$ws.Range("M132").Value = -2400.737
$ws.Range("N132").Value = -10545.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1288.9546
$ws.Range("I86").Value = 1321.0588
$ws.Range("K86").Value = 1321.0588
$ws.Range("M86").Value = -198.0588
$ws.Range("H89").Value = 1288.9546
$ws.Range("I89").Value = 1321.0588
$ws.Range("K89").Value = 6605.294
$ws.Range("M89").Value = -989.2939999999999
$ws.Range("H107").Value = 26601.7
$ws.Range("I107").Value = 27928.158
$ws.Range("J107").Value = 1399
$ws.Range("K107").Value = 27928.158
$ws.Range("L107").Value = 1399
$ws.Range("M107").Value = -26008.158
$ws.Range("N107").Value = -5239

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2033.2858
$ws.Range("I16").Value = 1746.8334
$ws.Range("J16").Value = 2415.2222
$ws.Range("K16").Value = 1746.8334
$ws.Range("L16").Value = 2415.2222
$ws.Range("M16").Value = -1459.8334
$ws.Range("N16").Value = -2989.2222
$ws.Range("H22").Value = 792.1111
$ws.Range("I22").Value = 924.75
$ws.Range("K22").Value = 924.75
$ws.Range("M22").Value = -574.75
$ws.Range("H62").Value = 10429.765
$ws.Range("I62").Value = 4756.4287
$ws.Range("K62").Value = 4756.4287
$ws.Range("M62").Value = -4132.4287
$ws.Range("H65").Value = 10429.765
$ws.Range("I65").Value = 4756.4287
$ws.Range("K65").Value = 23782.1435
$ws.Range("M65").Value = -20662.1435
$ws.Range("H88").Value = 32875
$ws.Range("J88").Value = 32875
$ws.Range("L88").Value = 32875
$ws.Range("N88").Value = -33687
$ws.Range("H91").Value = 32875
$ws.Range("J91").Value = 32875
$ws.Range("L91").Value = 32875
$ws.Range("N91").Value = -35683
$ws.Range("H107").Value = 320.81482
$ws.Range("I107").Value = 240.94737
$ws.Range("J107").Value = 510.5
$ws.Range("K107").Value = 240.94737
$ws.Range("L107").Value = 510.5
$ws.Range("M107").Value = 1679.05263
$ws.Range("N107").Value = -4350.5
$ws.Range("H113").Value = 2033.2858
$ws.Range("I113").Value = 1746.8334
$ws.Range("J113").Value = 2415.2222
$ws.Range("K113").Value = 1746.8334
$ws.Range("L113").Value = 2415.2222
$ws.Range("M113").Value = 423.1666
$ws.Range("N113").Value = -6755.2222
$ws.Range("H127").Value = 149999
$ws.Range("J127").Value = 149999
$ws.Range("L127").Value = 149999
$ws.Range("N127").Value = -159919
$ws.Range("H134").Value = 3231.9048
$ws.Range("I134").Value = 2305.923
$ws.Range("K134").Value = 6917.768999999999
$ws.Range("M134").Value = -4382.768999999999
$ws.Range("H137").Value = 88202
$ws.Range("J137").Value = 88202
$ws.Range("L137").Value = 88202
$ws.Range("N137").Value = -98402

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 2312.5652
$ws.Range("J113").Value = 1935.625
$ws.Range("L113").Value = 5806.875
$ws.Range("N113").Value = -10146.875
$ws.Range("H131").Value = 130278.45
$ws.Range("J131").Value = 1734
$ws.Range("L131").Value = 5202
$ws.Range("N131").Value = -15282

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5573.4136
$ws.Range("I70").Value = 5102.6665
$ws.Range("K70").Value = 5102.6665
$ws.Range("M70").Value = -4832.6665
$ws.Range("H73").Value = 5573.4136
$ws.Range("I73").Value = 5102.6665
$ws.Range("K73").Value = 5102.6665
$ws.Range("M73").Value = -4166.6665
$ws.Range("H80").Value = 6711.6
$ws.Range("I80").Value = 4079
$ws.Range("J80").Value = 11976.8
$ws.Range("K80").Value = 4079
$ws.Range("L80").Value = 11976.8
$ws.Range("M80").Value = -3081
$ws.Range("N80").Value = -13972.8
$ws.Range("H83").Value = 6711.6
$ws.Range("I83").Value = 4079
$ws.Range("J83").Value = 11976.8
$ws.Range("K83").Value = 20395
$ws.Range("L83").Value = 59884
$ws.Range("M83").Value = -15403
$ws.Range("N83").Value = -69868
$ws.Range("H97").Value = 1409.4103
$ws.Range("I97").Value = 1182.6
$ws.Range("J97").Value = 1814.4286
$ws.Range("K97").Value = 1182.6
$ws.Range("L97").Value = 1814.4286
$ws.Range("M97").Value = -686.5999999999999
$ws.Range("N97").Value = -2806.4286
$ws.Range("H126").Value = 2332.4285
$ws.Range("I126").Value = 2349.0908
$ws.Range("J126").Value = 2271.3333
$ws.Range("K126").Value = 7047.2724
$ws.Range("L126").Value = 6813.999899999999
$ws.Range("M126").Value = -4577.2724
$ws.Range("N126").Value = -11753.9999
$ws.Range("H128").Value = 103000
$ws.Range("J128").Value = 103000
$ws.Range("L128").Value = 103000
$ws.Range("N128").Value = -112960
$ws.Range("H136").Value = 8459.552
$ws.Range("J136").Value = 8459.552
$ws.Range("L136").Value = 25378.656
$ws.Range("N136").Value = -30478.656

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3117.0715
$ws.Range("I7").Value = 3134.2
$ws.Range("J7").Value = 3074.25
$ws.Range("K7").Value = 3134.2
$ws.Range("L7").Value = 3074.25
$ws.Range("M7").Value = -3022.2
$ws.Range("N7").Value = -3298.25
$ws.Range("H46").Value = 4980.6313
$ws.Range("I46").Value = 1106.75
$ws.Range("J46").Value = 7798
$ws.Range("K46").Value = 1106.75
$ws.Range("L46").Value = 7798
$ws.Range("M46").Value = -918.75
$ws.Range("N46").Value = -8174
$ws.Range("H55").Value = 1647.6957
$ws.Range("I55").Value = 393.54544
$ws.Range("J55").Value = 2797.3333
$ws.Range("K55").Value = 393.54544
$ws.Range("L55").Value = 2797.3333
$ws.Range("M55").Value = -220.54544
$ws.Range("N55").Value = -3143.3333
$ws.Range("H93").Value = 720121.1
$ws.Range("I93").Value = 1118914.6
$ws.Range("J93").Value = 2292.8
$ws.Range("K93").Value = 1118914.6
$ws.Range("L93").Value = 2292.8
$ws.Range("M93").Value = -1117666.6
$ws.Range("N93").Value = -4788.8
$ws.Range("H126").Value = 3117.0715
$ws.Range("I126").Value = 3134.2
$ws.Range("J126").Value = 3074.25
$ws.Range("K126").Value = 9402.599999999999
$ws.Range("L126").Value = 9222.75
$ws.Range("M126").Value = -6932.599999999999
$ws.Range("N126").Value = -14162.75
$ws.Range("H132").Value = 2220.3157
$ws.Range("I132").Value = 2339.423
$ws.Range("J132").Value = 1962.25
$ws.Range("K132").Value = 7018.268999999999
$ws.Range("L132").Value = 5886.75
$ws.Range("M132").Value = -4488.268999999999
$ws.Range("N132").Value = -10946.75
$ws.Range("H136").Value = 2728.55
$ws.Range("I136").Value = 2428.8076
$ws.Range("J136").Value = 3285.2144
$ws.Range("K136").Value = 7286.4228
$ws.Range("L136").Value = 9855.643199999999
$ws.Range("M136").Value = -4736.4228
$ws.Range("N136").Value = -14955.6432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 44430.484
$ws.Range("I122").Value = 74061.664
$ws.Range("K122").Value = 222184.992
$ws.Range("M122").Value = -219734.992
$ws.Range("H126").Value = 209876.62
$ws.Range("I126").Value = 1477
$ws.Range("K126").Value = 4431
$ws.Range("M126").Value = -1961
$ws.Range("H132").Value = 12185.357
$ws.Range("I132").Value = 14353.534
$ws.Range("J132").Value = 1705.8334
$ws.Range("K132").Value = 43060.602
$ws.Range("L132").Value = 5117.5002
$ws.Range("M132").Value = -40530.602
$ws.Range("N132").Value = -10177.5002
